# NIT-8020027654 Estado de Cuenta update
# - Update summary header values (Valor Mora total, Cant. Trabajadores, Cant. Periodos)
# - Replace the worker/period detail table (rows 16-77) with the refreshed dataset
# - Relocate the signature footer block from rows 76-77 down to rows 82-83

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Relocate the signature footer block (old rows 76-77) to rows 82-83 ---
# Copy formatting first, before the old row positions get reused by new data rows.
$ws.Range("B76:C77").Copy()
$ws.Range("B82:C83").PasteSpecial(-4122)
$ws.Range("H76:J77").Copy()
$ws.Range("H82:J83").PasteSpecial(-4122)

$ws.Range("B82").Value = "___________________________________"
$ws.Range("H82").Value = "___________________________________"
$ws.Range("B83").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H83").Value = "FIRMA DEL REPRESENTANTE LEGAL"

$ws.Range("B76:C76").UnMerge()
$ws.Range("H76:J76").UnMerge()
$ws.Range("B77:C77").UnMerge()
$ws.Range("H77:J77").UnMerge()

# Old rows 76-77 are now plain data rows again - drop the leftover footer text
$ws.Range("B76:J77").ClearContents()

$ws.Range("B82:C82").Merge()
$ws.Range("H82:J82").Merge()
$ws.Range("B83:C83").Merge()
$ws.Range("H83:J83").Merge()

# --- 2) Preserve the emphasised "last row" style (old row 71) by copying it to row 77 ---
# (row 77 becomes the new special last data row in the refreshed table)
$ws.Range("B71:J71").Copy()
$ws.Range("B77:J77").PasteSpecial(-4122)

# --- 3) Stamp the regular data-row style (from row 70) across rows 71-76 ---
# (row 71 stops being the special row; rows 72-76 are brand new rows)
$ws.Range("B70:J70").Copy()
$ws.Range("B71:J71").PasteSpecial(-4122)
$ws.Range("B72:J76").PasteSpecial(-4122)

# --- 4) Write the refreshed worker/period detail table into rows 16-77 ---
$tableData = @(
    @(16, "CC", "73187481", "ALEXANDER RAFAEL SANCHEZ DE HOYOS", "2304", 24640, 781242),
    @(17, "CC", "73187481", "ALEXANDER RAFAEL SANCHEZ DE HOYOS", "2303", 24640, 781242),
    @(18, "CC", "1047418375", "LUIS GABRIEL ARNEDO YANES", "2401", 46400, 1160000),
    @(19, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "2003", 31249, 781242),
    @(20, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "2002", 31249, 781242),
    @(21, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "2001", 31249, 781242),
    @(22, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1912", 31249, 781242),
    @(23, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1911", 31249, 781242),
    @(24, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1910", 31249, 781242),
    @(25, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1909", 31249, 781242),
    @(26, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1908", 31249, 781242),
    @(27, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1907", 31249, 781242),
    @(28, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1906", 31249, 781242),
    @(29, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1905", 31249, 781242),
    @(30, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1904", 31249, 781242),
    @(31, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1903", 31249, 781242),
    @(32, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1902", 31249, 781242),
    @(33, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1901", 31249, 781242),
    @(34, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1812", 31249, 781242),
    @(35, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1811", 31249, 781242),
    @(36, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1810", 31249, 781242),
    @(37, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1809", 31249, 781242),
    @(38, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1808", 24640, 781242),
    @(39, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1807", 24640, 781242),
    @(40, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1806", 24640, 781242),
    @(41, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1805", 24640, 781242),
    @(42, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1804", 24640, 781242),
    @(43, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1803", 24640, 781242),
    @(44, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1802", 24640, 781242),
    @(45, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1801", 24640, 781242),
    @(46, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1712", 24640, 781242),
    @(47, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1711", 24640, 781242),
    @(48, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1710", 24640, 781242),
    @(49, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1709", 24640, 781242),
    @(50, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1708", 24640, 781242),
    @(51, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1707", 24640, 781242),
    @(52, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1706", 24640, 781242),
    @(53, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1705", 24640, 781242),
    @(54, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1704", 24640, 781242),
    @(55, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1703", 24640, 781242),
    @(56, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1702", 24640, 781242),
    @(57, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1701", 24640, 781242),
    @(58, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1612", 24640, 781242),
    @(59, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1611", 24640, 781242),
    @(60, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1610", 24640, 781242),
    @(61, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1609", 24640, 781242),
    @(62, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1608", 24640, 781242),
    @(63, "CC", "12685506", "HERNAN YESID MONTERO JEJEN", "1607", 24640, 781242),
    @(64, "CC", "1049829313", "DAICER DANIEL PRENS ARIAS", "2103", 36341, 908526),
    @(65, "CC", "1049829313", "DAICER DANIEL PRENS ARIAS", "2102", 36341, 908526),
    @(66, "CC", "1143336332", "GILLIANO MARCO SERNA HERNANDEZ", "2110", 18170, 908526),
    @(67, "CC", "1049537251", "CARLOS ANDRES MARTINEZ TABORDA", "2311", 40000, 1000000),
    @(68, "CC", "1049537251", "CARLOS ANDRES MARTINEZ TABORDA", "2310", 40000, 1000000),
    @(69, "CC", "1049537251", "CARLOS ANDRES MARTINEZ TABORDA", "2309", 40000, 1000000),
    @(70, "CC", "1049537251", "CARLOS ANDRES MARTINEZ TABORDA", "2308", 40000, 1000000),
    @(71, "CC", "1049537251", "CARLOS ANDRES MARTINEZ TABORDA", "2307", 40000, 1000000),
    @(72, "CC", "1049537251", "CARLOS ANDRES MARTINEZ TABORDA", "2306", 40000, 1000000),
    @(73, "CC", "73270377", "OSCAR MANUEL HERNANDEZ VARGAS", "1905", 48000, 1200000),
    @(74, "CC", "73270377", "OSCAR MANUEL HERNANDEZ VARGAS", "1904", 48000, 1200000),
    @(75, "CC", "73270377", "OSCAR MANUEL HERNANDEZ VARGAS", "1903", 48000, 1200000),
    @(76, "CC", "73270377", "OSCAR MANUEL HERNANDEZ VARGAS", "1902", 48000, 1200000),
    @(77, "PPT", "1239399", "ANTHONY YEFERSSON VILLANUEVA GIL", "2306", 40000, 1300000)
)
foreach ($r in $tableData) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 2).Value = $r[1]
    $ws.Cells.Item($rowNum, 3).Value = $r[2]
    $ws.Cells.Item($rowNum, 4).Value = $r[3]
    $ws.Cells.Item($rowNum, 5).Value = $r[4]
    $ws.Cells.Item($rowNum, 6).Value = $r[5]
    $ws.Cells.Item($rowNum, 7).Value = $r[6]
}


# --- 5) Update the summary figures above the table ---
$ws.Range("E11").Value = 1936423
$ws.Range("C13").Value = 8
$ws.Range("F13").Value = 57
